$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4351974890530812
$ws.Range("D2").Value = 0.1560007319877279
$ws.Range("E2").Value = 0.1217164210942254
$ws.Range("F2").Value = 1.235639417020593
$ws.Range("G2").Value = 0.00242339283234936
$ws.Range("I2").Value = 0.9925639220564761
$ws.Range("J2").Value = 0.1296473592992839
$ws.Range("M2").Value = 1.885147280768791
$ws.Range("O2").Value = 3.122355927703978
$ws.Range("B3").Value = 0.3820937015902643
$ws.Range("D3").Value = 0.1580539028583967
$ws.Range("E3").Value = 0.1240325594675937
$ws.Range("F3").Value = 1.231589265483009
$ws.Range("G3").Value = 0.002427497315677646
$ws.Range("I3").Value = 0.9094617385134427
$ws.Range("J3").Value = 0.1321328220062661
$ws.Range("M3").Value = 1.67582380016789
$ws.Range("O3").Value = 3.085961469642569
$ws.Range("B4").Value = 0.3493570768498273
$ws.Range("D4").Value = 0.1593854234050909
$ws.Range("E4").Value = 0.1255513087736198
$ws.Range("F4").Value = 1.230175737371752
$ws.Range("G4").Value = 0.00243015039618329
$ws.Range("I4").Value = 0.8588165632416036
$ws.Range("J4").Value = 0.1337940835872544
$ws.Range("M4").Value = 1.546669850851501
$ws.Range("O4").Value = 3.066348599718197
$ws.Range("B5").Value = 0.3359844396274525
$ws.Range("D5").Value = 0.1599458573185792
$ws.Range("E5").Value = 0.1261944831075912
$ws.Range("F5").Value = 1.229868270388323
$ws.Range("G5").Value = 0.002431265077957979
$ws.Range("I5").Value = 0.838275356550497
$ws.Range("J5").Value = 0.1345049754737033
$ws.Range("M5").Value = 1.493885221971723
$ws.Range("O5").Value = 3.059039796197055
$ws.Range("B6").Value = 0.3337619988598703
$ws.Range("D6").Value = 0.1600399939971631
$ws.Range("E6").Value = 0.1263027468860329
$ws.Range("F6").Value = 1.229833397127351
$ws.Range("G6").Value = 0.00243145219810731
$ws.Range("I6").Value = 0.834870418258717
$ws.Range("J6").Value = 0.134625064646702
$ws.Range("M6").Value = 1.485111252359175
$ws.Range("O6").Value = 3.057867353128927
$ws.Range("B7").Value = 0.3491768582516386
$ws.Range("D7").Value = 0.1593929094139312
$ws.Range("E7").Value = 0.1255598846097641
$ws.Range("F7").Value = 1.230170505154959
$ws.Range("G7").Value = 0.002430165293257044
$ws.Range("I7").Value = 0.8585391419827033
$ws.Range("J7").Value = 0.1338035337102284
$ws.Range("M7").Value = 1.545958593504167
$ws.Range("O7").Value = 3.066247267759508
$ws.Range("B8").Value = 0.4169148022052411
$ws.Range("D8").Value = 0.1566939627412474
$ws.Range("E8").Value = 0.1224949536183892
$ws.Range("F8").Value = 1.234019360057573
$ws.Range("G8").Value = 0.002424780537679588
$ws.Range("I8").Value = 0.9638323407655918
$ws.Range("J8").Value = 0.1304762438869886
$ws.Range("M8").Value = 1.81310571497815
$ws.Range("O8").Value = 3.109237503146943
$ws.Range("B9").Value = 0.5486900375022685
$ws.Range("D9").Value = 0.1519631498232474
$ws.Range("E9").Value = 0.117252542249819
$ws.Range("F9").Value = 1.250145241934092
$ws.Range("G9").Value = 0.002415270686076605
$ws.Range("I9").Value = 1.173271505546751
$ws.Range("J9").Value = 0.1250275630467357
$ws.Range("M9").Value = 2.331815097552919
$ws.Range("O9").Value = 3.215408698877297
$ws.Range("B10").Value = 0.6448416710211404
$ws.Range("D10").Value = 0.1488291612563959
$ws.Range("E10").Value = 0.1138706766853748
$ws.Range("F10").Value = 1.267312135969036
$ws.Range("G10").Value = 0.002408916730739312
$ws.Range("I10").Value = 1.328898680605164
$ws.Range("J10").Value = 0.121685471592393
$ws.Range("M10").Value = 2.70956248262894
$ws.Range("O10").Value = 3.307004232304052
$ws.Range("B11").Value = 0.6884366492348022
$ws.Range("D11").Value = 0.1474774688590266
$ws.Range("E11").Value = 0.1124345339683241
$ws.Range("F11").Value = 1.276296414246247
$ws.Range("G11").Value = 0.002406162102392108
$ws.Range("I11").Value = 1.400069516993454
$ws.Range("J11").Value = 0.1203097344800739
$ws.Range("M11").Value = 2.8806439028169
$ws.Range("O11").Value = 3.351681870514085
$ws.Range("B12").Value = 0.7049236853720799
$ws.Range("D12").Value = 0.1469762460676485
$ws.Range("E12").Value = 0.1119054465612255
$ws.Range("F12").Value = 1.279869022042405
$ws.Range("G12").Value = 0.002405138413651761
$ws.Range("I12").Value = 1.427073043764551
$ws.Range("J12").Value = 0.1198096688726906
$ws.Range("M12").Value = 2.945315102675835
$ws.Range("O12").Value = 3.369037422755241
$ws.Range("B13").Value = 0.7013738698407792
$ws.Range("D13").Value = 0.1470837204814597
$ws.Range("E13").Value = 0.1120187383999607
$ws.Range("F13").Value = 1.279091992009469
$ws.Range("G13").Value = 0.002405358020806316
$ws.Range("I13").Value = 1.421255028041799
$ws.Range("J13").Value = 0.1199164359166076
$ws.Range("M13").Value = 2.931392129230318
$ws.Range("O13").Value = 3.365280086557107
$ws.Range("B14").Value = 0.6897934805895716
$ws.Range("D14").Value = 0.1474360198614502
$ws.Range("E14").Value = 0.1123907098102688
$ws.Range("F14").Value = 1.27658691003532
$ws.Range("G14").Value = 0.002406077494166168
$ws.Range("I14").Value = 1.402290064796148
$ws.Range("J14").Value = 0.1202681743914411
$ws.Range("M14").Value = 2.885966744854215
$ws.Range("O14").Value = 3.35310093788911
$ws.Range("B15").Value = 0.6826973487570172
$ws.Range("D15").Value = 0.1476531981285945
$ws.Range("E15").Value = 0.1126204749596589
$ws.Range("F15").Value = 1.275074718027057
$ws.Range("G15").Value = 0.002406520719511523
$ws.Range("I15").Value = 1.390680300234322
$ws.Range("J15").Value = 0.1204863487662493
$ws.Range("M15").Value = 2.858127451743428
$ws.Range("O15").Value = 3.345697906448038
$ws.Range("B16").Value = 0.6419896692275699
$ws.Range("D16").Value = 0.1489189857772555
$ws.Range("E16").Value = 0.1139665936429566
$ws.Range("F16").Value = 1.2667487635036
$ws.Range("G16").Value = 0.002409099476505016
$ws.Range("I16").Value = 1.324254938575422
$ws.Range("J16").Value = 0.1217782976907742
$ws.Range("M16").Value = 2.698366261453174
$ws.Range("O16").Value = 3.304145379504803
$ws.Range("B17").Value = 0.6169792821976898
$ws.Range("D17").Value = 0.1497144495604346
$ws.Range("E17").Value = 0.1148186198288421
$ws.Range("F17").Value = 1.261943061508362
$ws.Range("G17").Value = 0.002410716174514034
$ws.Range("I17").Value = 1.283600304264809
$ws.Range("J17").Value = 0.1226079740847936
$ws.Range("M17").Value = 2.600160525546784
$ws.Range("O17").Value = 3.279428268246818
$ws.Range("B18").Value = 0.602580366464764
$ws.Range("D18").Value = 0.1501789423604798
$ws.Range("E18").Value = 0.1153183078402353
$ws.Range("F18").Value = 1.259289452860273
$ws.Range("G18").Value = 0.002411658846273048
$ws.Range("I18").Value = 1.260252226694902
$ws.Range("J18").Value = 0.123098785851667
$ws.Range("M18").Value = 2.54360413261864
$ws.Range("O18").Value = 3.265494736748792
$ws.Range("B19").Value = 0.5977028235047328
$ws.Range("D19").Value = 0.1503374076683208
$ws.Range("E19").Value = 0.1154891455187794
$ws.Range("F19").Value = 1.258409920230534
$ws.Range("G19").Value = 0.002411980218557659
$ws.Range("I19").Value = 1.252353087927347
$ws.Range("J19").Value = 0.1232672996434303
$ws.Range("M19").Value = 2.524443049930682
$ws.Range("O19").Value = 3.26082557659123
$ws.Range("B20").Value = 0.6196430928468715
$ws.Range("D20").Value = 0.1496290504609075
$ws.Range("E20").Value = 0.1147269237424009
$ws.Range("F20").Value = 1.262443189150886
$ws.Range("G20").Value = 0.002410542751275892
$ws.Range("I20").Value = 1.287924406012252
$ws.Range("J20").Value = 0.1225182449339286
$ws.Range("M20").Value = 2.610622081279246
$ws.Range("O20").Value = 3.282030117089164
$ws.Range("B21").Value = 0.6931955080657417
$ws.Range("D21").Value = 0.1473322524256417
$ws.Range("E21").Value = 0.1122810521908253
$ws.Range("F21").Value = 1.277318075224571
$ws.Range("G21").Value = 0.002405865640789572
$ws.Range("I21").Value = 1.407859114497199
$ws.Range("J21").Value = 0.1201642923244037
$ws.Range("M21").Value = 2.899312401701195
$ws.Range("O21").Value = 3.356666352612194
$ws.Range("B22").Value = 0.7411410562544916
$ws.Range("D22").Value = 0.1458931422867238
$ws.Range("E22").Value = 0.1107685195383361
$ws.Range("F22").Value = 1.288033924653149
$ws.Range("G22").Value = 0.002402922081295941
$ws.Range("I22").Value = 1.486549868865495
$ws.Range("J22").Value = 0.1187476966377687
$ws.Range("M22").Value = 3.087323843521546
$ws.Range("O22").Value = 3.407995174022005
$ws.Range("B23").Value = 0.7155632637869189
$ws.Range("D23").Value = 0.1466555524732396
$ws.Range("E23").Value = 0.111567905472155
$ws.Range("F23").Value = 1.282223203531174
$ws.Range("G23").Value = 0.002404482789654791
$ws.Range("I23").Value = 1.444523494264615
$ws.Range("J23").Value = 0.1194925765825516
$ws.Range("M23").Value = 2.987040905173757
$ws.Range("O23").Value = 3.380365317808184
$ws.Range("B24").Value = 0.6184388466499797
$ws.Range("D24").Value = 0.149667637059423
$ws.Range("E24").Value = 0.1147683488735183
$ws.Range("F24").Value = 1.262216741242014
$ws.Range("G24").Value = 0.00241062111480133
$ws.Range("I24").Value = 1.285969402226129
$ws.Range("J24").Value = 0.1225587684422962
$ws.Range("M24").Value = 2.605892712495148
$ws.Range("O24").Value = 3.280852959759386
$ws.Range("B25").Value = 0.5131566668782455
$ws.Range("D25").Value = 0.1531828924361758
$ws.Range("E25").Value = 0.1185883816168491
$ws.Range("F25").Value = 1.244855155769315
$ws.Range("G25").Value = 0.002417731702870642
$ws.Range("I25").Value = 1.116301487558587
$ws.Range("J25").Value = 0.1263858584964233
$ws.Range("M25").Value = 2.192064399784016
$ws.Range("O25").Value = 3.18431865987867
